$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6844159960746765
$ws.Range("B1").Value = 1.560806274414062
$ws.Range("C1").Value = 7.146332740783691
$ws.Range("D1").Value = 2.819880247116089
$ws.Range("E1").Value = 1.647903919219971
